# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Feria Lagunitas de Puerto Montt - Limón"
# right after the existing row 406, pushing the subsequent rows down by two
# positions (old row 407 -> new row 409, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 407-408; everything from the old row 407 onward
# shifts down to 409 onward.
$ws.Rows("407:408").Insert()

# New row 407
$ws.Range("A407").Value = 4
$ws.Range("B407").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C407").Value = "Los Lagos"
$ws.Range("D407").Value = 44694
$ws.Range("E407").Value = 10
$ws.Range("F407").Value = "Fruta"
$ws.Range("G407").Value = 100102
$ws.Range("H407").Value = "Cítricos"
$ws.Range("I407").Value = 100102003
$ws.Range("J407").Value = "Limón"
$ws.Range("K407").Value = "Sin especificar"
$ws.Range("L407").Value = "1a plateado"
$ws.Range("M407").Value = 1000
$ws.Range("N407").Value = 24000
$ws.Range("O407").Value = 25000
$ws.Range("P407").Value = 24500
$ws.Range("Q407").Value = "`$/malla 18 kilos"
$ws.Range("R407").Value = "Provincia de Melipilla"
$ws.Range("S407").Value = 1361
$ws.Range("T407").Value = 18

# New row 408
$ws.Range("A408").Value = 4
$ws.Range("B408").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C408").Value = "Los Lagos"
$ws.Range("D408").Value = 44694
$ws.Range("E408").Value = 10
$ws.Range("F408").Value = "Fruta"
$ws.Range("G408").Value = 100102
$ws.Range("H408").Value = "Cítricos"
$ws.Range("I408").Value = 100102003
$ws.Range("J408").Value = "Limón"
$ws.Range("K408").Value = "Sin especificar"
$ws.Range("L408").Value = "2a plateado"
$ws.Range("M408").Value = 500
$ws.Range("N408").Value = 21000
$ws.Range("O408").Value = 21000
$ws.Range("P408").Value = 21000
$ws.Range("Q408").Value = "`$/malla 18 kilos"
$ws.Range("R408").Value = "Provincia de Melipilla"
$ws.Range("S408").Value = 1167
$ws.Range("T408").Value = 18
